$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "color" column header and populate the unique values first, in the
# same order they should appear in the shared-strings table, then fill in
# the repeated values.
$ws.Range("H1").Value = "color"
$ws.Range("H2").Value = "blue"
$ws.Range("H3").Value = "orange"
$ws.Range("H5").Value = "green"
$ws.Range("H8").Value = "brown"
$ws.Range("H9").Value = "red"
$ws.Range("H6").Value = "yellow"
$ws.Range("H4").Value = "blue"
$ws.Range("H7").Value = "green"

$ws.Range("H11").Select()
